# Generate Report for Handoff
# Adds two new localization-status rows (a66573f6... and afba94ed...) to the
# Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/aa126b143a9a5b945044e8dfb8c614b80be484d0/e2e/"

# Excel auto-detects "True"/"False" typed into a cell as booleans and empty
# strings as "no value" - prefix with a leading apostrophe to force them to
# be stored as literal text, matching the source data.
$sTrue = "'True"
$sFalse = "'False"
$sEmpty = "'"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4 - a66573f6-afa9-4d68-993f-fc618f0c07e8.md
$wsOverview.Range("A4").Value = "a66573f6-afa9-4d68-993f-fc618f0c07e8.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), ($repoBase + "a66573f6-afa9-4d68-993f-fc618f0c07e8.md"), "", "", "e2e\a66573f6-afa9-4d68-993f-fc618f0c07e8.md") | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = $sEmpty
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2017-02-21 02:40:49"
$wsOverview.Range("G4").NumberFormat = $wsOverview.Range("G2").NumberFormat

# Row 5 - afba94ed-dc6f-42fa-a9df-84cd86f730c1.md
$wsOverview.Range("A5").Value = "afba94ed-dc6f-42fa-a9df-84cd86f730c1.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), ($repoBase + "afba94ed-dc6f-42fa-a9df-84cd86f730c1.md"), "", "", "e2e\afba94ed-dc6f-42fa-a9df-84cd86f730c1.md") | Out-Null
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = $sEmpty
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2017-02-21 02:40:49"
$wsOverview.Range("G5").NumberFormat = $wsOverview.Range("G2").NumberFormat

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

# Row 4 - a66573f6-afa9-4d68-993f-fc618f0c07e8.md
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($repoBase + "a66573f6-afa9-4d68-993f-fc618f0c07e8.md"), "", "", "a66573f6-afa9-4d68-993f-fc618f0c07e8.md") | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = $sFalse
$wsZh.Range("G4").Value = "a66573f6-afa9-4d68-993f-fc618f0c07e8.aa695096d6c0fc8bb16c7113fc9f24272b2762a1.zh-cn.xlf"
$wsZh.Range("H4").Value = "2017-02-21 02:40:34"
$wsZh.Range("H4").NumberFormat = $wsZh.Range("H2").NumberFormat
$wsZh.Range("I4").Value = $sEmpty
$wsZh.Range("J4").Value = $sEmpty
$wsZh.Range("K4").Value = $sEmpty
$wsZh.Range("L4").Value = "0001-01-01 00:00:00"
$wsZh.Range("L4").NumberFormat = $wsZh.Range("L2").NumberFormat
$wsZh.Range("M4").Value = $sEmpty
$wsZh.Range("N4").Value = $sEmpty
$wsZh.Range("O4").Value = $sTrue
$wsZh.Range("P4").Value = $sEmpty
$wsZh.Range("Q4").Value = $sFalse
$wsZh.Range("R4").Value = $sEmpty

# Row 5 - afba94ed-dc6f-42fa-a9df-84cd86f730c1.md
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($repoBase + "afba94ed-dc6f-42fa-a9df-84cd86f730c1.md"), "", "", "afba94ed-dc6f-42fa-a9df-84cd86f730c1.md") | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = $sFalse
$wsZh.Range("G5").Value = "afba94ed-dc6f-42fa-a9df-84cd86f730c1.676be2bd1953dad2ea6c8a03bd11e64f568ab728.zh-cn.xlf"
$wsZh.Range("H5").Value = "2017-02-21 02:40:34"
$wsZh.Range("H5").NumberFormat = $wsZh.Range("H2").NumberFormat
$wsZh.Range("I5").Value = $sEmpty
$wsZh.Range("J5").Value = $sEmpty
$wsZh.Range("K5").Value = $sEmpty
$wsZh.Range("L5").Value = "0001-01-01 00:00:00"
$wsZh.Range("L5").NumberFormat = $wsZh.Range("L2").NumberFormat
$wsZh.Range("M5").Value = $sEmpty
$wsZh.Range("N5").Value = $sEmpty
$wsZh.Range("O5").Value = $sTrue
$wsZh.Range("P5").Value = $sEmpty
$wsZh.Range("Q5").Value = $sFalse
$wsZh.Range("R5").Value = $sEmpty

$wsZh.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

# Row 4 - a66573f6-afa9-4d68-993f-fc618f0c07e8.md
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($repoBase + "a66573f6-afa9-4d68-993f-fc618f0c07e8.md"), "", "", "a66573f6-afa9-4d68-993f-fc618f0c07e8.md") | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = $sFalse
$wsDe.Range("G4").Value = "a66573f6-afa9-4d68-993f-fc618f0c07e8.aa695096d6c0fc8bb16c7113fc9f24272b2762a1.de-de.xlf"
$wsDe.Range("H4").Value = "2017-02-21 02:40:49"
$wsDe.Range("H4").NumberFormat = $wsDe.Range("H2").NumberFormat
$wsDe.Range("I4").Value = $sEmpty
$wsDe.Range("J4").Value = $sEmpty
$wsDe.Range("K4").Value = $sEmpty
$wsDe.Range("L4").Value = "0001-01-01 00:00:00"
$wsDe.Range("L4").NumberFormat = $wsDe.Range("L2").NumberFormat
$wsDe.Range("M4").Value = $sEmpty
$wsDe.Range("N4").Value = $sEmpty
$wsDe.Range("O4").Value = $sTrue
$wsDe.Range("P4").Value = $sEmpty
$wsDe.Range("Q4").Value = $sFalse
$wsDe.Range("R4").Value = $sEmpty

# Row 5 - afba94ed-dc6f-42fa-a9df-84cd86f730c1.md
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($repoBase + "afba94ed-dc6f-42fa-a9df-84cd86f730c1.md"), "", "", "afba94ed-dc6f-42fa-a9df-84cd86f730c1.md") | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = $sFalse
$wsDe.Range("G5").Value = "afba94ed-dc6f-42fa-a9df-84cd86f730c1.676be2bd1953dad2ea6c8a03bd11e64f568ab728.de-de.xlf"
$wsDe.Range("H5").Value = "2017-02-21 02:40:49"
$wsDe.Range("H5").NumberFormat = $wsDe.Range("H2").NumberFormat
$wsDe.Range("I5").Value = $sEmpty
$wsDe.Range("J5").Value = $sEmpty
$wsDe.Range("K5").Value = $sEmpty
$wsDe.Range("L5").Value = "0001-01-01 00:00:00"
$wsDe.Range("L5").NumberFormat = $wsDe.Range("L2").NumberFormat
$wsDe.Range("M5").Value = $sEmpty
$wsDe.Range("N5").Value = $sEmpty
$wsDe.Range("O5").Value = $sTrue
$wsDe.Range("P5").Value = $sEmpty
$wsDe.Range("Q5").Value = $sFalse
$wsDe.Range("R5").Value = $sEmpty

$wsDe.Columns.Item(3).ColumnWidth = 16.3
